$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("G2")
$v = $c.Value2
Write-Host "V2: $v"
$v3 = $c.Text
Write-Host "TEXT: $v3"
